$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Test Case #Id" column (A) was renumbered into one continuous
# sequence (1..19) for rows 5-23, filling in rows that previously had
# no number at all.
for ($i = 0; $i -lt 19; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}

# Row 23's first cell (A23) used a stray/duplicate cell style (identical
# in appearance to the style used by every other numbered row, A5:A22,
# just carrying a redundant fill flag). Normalize it so it matches the
# rest of the column exactly.
$a23 = $ws.Range("A23")
$a23.HorizontalAlignment = -4108   # xlCenter
$a23.WrapText = $true
$a23.Font.Size = 16
